# "Colocando header nos graficos"
# Adds a header label to column A (row 1) on each data sheet, removes the
# bold/border header styling that had incorrectly been applied to the
# row-label cells (A2:A12 / A2:A3 / A2:A3), fixes a few accented Portuguese
# words, and (on the "Emissoes Totais" sheet) drops the unused "Teto" row,
# and (on the "Custo Total" sheet) relabels/updates the cost-by-expansion
# figures.

$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: "Potencia Acumulada", "Geracao Periodo Medio",
#     "Atendimento a Ponta", "Potencia Incremental" share the same layout:
#     column A holds the technology/source name, rows 2-12.
$sourceSheets = @(1, 2, 3, 4)
foreach ($idx in $sourceSheets) {
    $ws = $wb.Worksheets.Item($idx)

    # Add the column header in A1, matching the styling already used by
    # the other header cells (B1:E1) on the same row.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # The row labels (A2:A12) should not carry the bold/bordered header
    # style - strip it back to the default.
    $ws.Range("A2:A12").ClearFormats()

    # Fix missing accents.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

# --- Sheet 5: "Emissoes Totais (MtCO2eq)"
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2:A3").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

# The "Teto" row is no longer used - remove it entirely.
$ws5.Rows("4:4").Delete()

# --- Sheet 6: "Custo Total (bilhões de R$)"
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's new label ("2015") looks numeric, so force text formatting before
# assigning it, then re-apply the original bold/bordered header style
# (NumberFormat changes spin up a distinct style record).
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
$ws6.Range("A2").Copy()
$ws6.Range("B1").PasteSpecial(-4122)

$ws6.Range("A2:A3").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 602
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

Write-Host "done"
